$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final grid values for rows 2-7 (columns A-G), reflecting a "delete row" on the
# inventory grid followed by a refresh of the bound data (commit: "added delete
# on all grids"). Column F = Date Added, Column G = Last Updated.
$rows = @(
    @{ Row=2; A=2; B=2;  C=3; D=1; E=40;   F="25/03/2023"; G="25/03/2023" },
    @{ Row=3; A=3; B=10; C=3; D=1; E=899;  F="25/03/2023"; G="25/03/2023" },
    @{ Row=4; A=4; B=1;  C=1; D=1; E=8000; F="01/01/0001"; G="25/03/2023" },
    @{ Row=5; A=5; B=1;  C=1; D=1; E=51;   F="25/03/2023"; G="25/03/2023" },
    @{ Row=6; A=8; B=8;  C=1; D=1; E=1;    F="01/01/0001"; G="25/03/2023" },
    @{ Row=7; A=1; B=10; C=2; D=3; E=15;   F="01/01/0001"; G="25/03/2023" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E

    # Columns F/G hold plain text "dates" (shared strings, not real date
    # serials). "01/01/0001" parses as a valid date otherwise, so force a
    # text format while assigning it, then restore the default style so the
    # cell still comes out with the workbook's normal (style 0) formatting.
    $fCell = $ws.Cells.Item($row, 6)
    $fCell.NumberFormat = "@"
    $fCell.Value = $r.F
    $fCell.Style = "Normal"

    $gCell = $ws.Cells.Item($row, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = $r.G
    $gCell.Style = "Normal"
}
